$d = $word.ActiveDocument

$replacements = @(
    @('75×56=4200', '23×73=1679'),
    @('17×34=578', '69×85=5865'),
    @('56×62=3472', '32×78=2496'),
    @('66×85=5610', '18×42=756'),
    @('25×15=375', '11×30=330'),
    @('55×79=4345', '99×20=1980'),
    @('50×36=1800', '42×25=1050'),
    @('54×13=702', '59×26=1534'),
    @('53×48=2544', '82×78=6396'),
    @('37×19=703', '19×97=1843'),
    @('56×99=5544', '18×75=1350'),
    @('75×86=6450', '44×39=1716'),
    @('97×82=7954', '57×18=1026'),
    @('93×20=1860', '30×72=2160'),
    @('42×17=714', '52×21=1092'),
    @('43×36=1548', '19×67=1273'),
    @('17×90=1530', '95×17=1615'),
    @('52×81=4212', '60×99=5940'),
    @('70×73=5110', '99×93=9207'),
    @('72×30=2160', '53×46=2438'),
    @('47×17=799', '24×22=528'),
    @('40×93=3720', '82×12=984'),
    @('48×68=3264', '14×31=434'),
    @('46×90=4140', '39×51=1989'),
    @('45×61=2745', '96×24=2304'),
    @('92×41=3772', '82×75=6150'),
    @('14×52=728', '27×28=756'),
    @('32×19=608', '83×94=7802'),
    @('98×78=7644', '62×49=3038'),
    @('72×74=5328', '61×34=2074'),
    @('13×52=676', '70×16=1120'),
    @('44×19=836', '97×28=2716'),
    @('38×46=1748', '84×45=3780'),
    @('90×27=2430', '75×76=5700'),
    @('56×80=4480', '72×41=2952'),
    @('60×41=2460', '50×92=4600'),
    @('90×70=6300', '77×74=5698'),
    @('19×92=1748', '20×40=800'),
    @('60×68=4080', '73×12=876'),
    @('54×74=3996', '71×89=6319'),
    @('50×39=1950', '71×93=6603'),
    @('70×63=4410', '44×34=1496'),
    @('19×61=1159', '22×14=308'),
    @('22×75=1650', '54×17=918'),
    @('99×37=3663', '75×28=2100'),
    @('71×54=3834', '72×64=4608'),
    @('15×11=165', '25×55=1375'),
    @('25×50=1250', '72×67=4824'),
    @('84×75=6300', '12×49=588'),
    @('47×95=4465', '60×13=780'),
    @('57×84=4788', '26×77=2002'),
    @('16×50=800', '100×82=8200'),
    @('98×82=8036', '12×20=240'),
    @('94×53=4982', '81×44=3564'),
    @('55×20=1100', '87×74=6438'),
    @('70×35=2450', '39×92=3588'),
    @('25×71=1775', '62×74=4588'),
    @('36×41=1476', '21×100=2100'),
    @('88×27=2376', '38×44=1672'),
    @('76×30=2280', '47×74=3478'),
    @('67×53=3551', '69×55=3795'),
    @('50×43=2150', '50×40=2000'),
    @('94×91=8554', '43×85=3655'),
    @('62×21=1302', '92×24=2208'),
    @('83×66=5478', '37×84=3108'),
    @('70×50=3500', '55×50=2750'),
    @('79×15=1185', '62×13=806'),
    @('58×16=928', '45×16=720'),
    @('48×60=2880', '57×13=741'),
    @('55×72=3960', '36×86=3096'),
    @('62×52=3224', '65×48=3120'),
    @('10×34=340', '70×51=3570'),
    @('69×66=4554', '65×43=2795'),
    @('53×79=4187', '20×52=1040'),
    @('16×93=1488', '73×80=5840'),
    @('90×39=3510', '19×59=1121'),
    @('78×39=3042', '45×73=3285'),
    @('24×82=1968', '22×38=836'),
    @('10×64=640', '14×71=994'),
    @('14×38=532', '63×46=2898'),
    @('45×89=4005', '99×82=8118'),
    @('87×40=3480', '97×34=3298'),
    @('40×25=1000', '12×36=432'),
    @('46×82=3772', '53×95=5035'),
    @('63×41=2583', '70×94=6580'),
    @('91×55=5005', '22×19=418'),
    @('32×23=736', '55×95=5225'),
    @('40×17=680', '97×45=4365'),
    @('55×46=2530', '88×89=7832'),
    @('51×89=4539', '32×71=2272'),
    @('49×30=1470', '31×90=2790'),
    @('66×18=1188', '39×59=2301'),
    @('41×47=1927', '96×61=5856'),
    @('59×64=3776', '69×48=3312'),
    @('76×15=1140', '83×81=6723'),
    @('51×53=2703', '11×61=671'),
    @('47×76=3572', '70×31=2170'),
    @('44×48=2112', '23×77=1771'),
    @('74×93=6882', '21×59=1239'),
    @('34×86=2924', '33×85=2805')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done: $($replacements.Count) replacements attempted"
